$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 22: bus stop info
$ws.Range("A22").Value = "105, 106 Bus Stop, Lancaster/Bryn Mawr"
$ws.Range("B22").Value = 94

$ws.Range("D22").Value = "https://public.tableau.com/app/profile/daniel.sandiford4261/viz/Stop_Summary_Report_Public_1/Ridership_Map_Public?publish=yes"
$ws.Hyperlinks.Add($ws.Range("D22"), "https://public.tableau.com/app/profile/daniel.sandiford4261/viz/Stop_Summary_Report_Public_1/Ridership_Map_Public?publish=yes", "", "", "https://public.tableau.com/app/profile/daniel.sandiford4261/viz/Stop_Summary_Report_Public_1/Ridership_Map_Public?publish=yes")
$ws.Range("D22").Style = "Hyperlink"

# Column widths
$ws.Columns.Item(1).ColumnWidth = 36.7109375
$ws.Columns.Item(4).ColumnWidth = 113.42578125

# View state
$ws.Range("D22").Select()
